$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3071
$ws.Range("I125").Value = 3097.5
$ws.Range("K125").Value = 27877.5
$ws.Range("M125").Value = -25417.5
$ws.Range("H133").Value = 43749.25
$ws.Range("J133").Value = 43749.25
$ws.Range("L133").Value = 43749.25
$ws.Range("N133").Value = -53869.25
$ws.Range("H137").Value = 1283.8667
$ws.Range("I137").Value = 1164.7059
$ws.Range("J137").Value = 1439.6923
$ws.Range("K137").Value = 3494.1177
$ws.Range("L137").Value = 4319.0769
$ws.Range("M137").Value = -944.1176999999998
$ws.Range("N137").Value = -9419.0769
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1967
$ws.Range("I32").Value = 1996.1837
$ws.Range("J32").Value = 537
$ws.Range("K32").Value = 1996.1837
$ws.Range("L32").Value = 537
$ws.Range("M32").Value = -1709.1837
$ws.Range("N32").Value = -1111
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H99").Value = 62501440
$ws.Range("I99").Value = 100001280
$ws.Range("K99").Value = 100001280
$ws.Range("M99").Value = -99999782
$ws.Range("H140").Value = 22277.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 22277.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 22277.5
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -32637.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 259
$ws.Range("I7").Value = 321.25
$ws.Range("J7").Value = 209.2
$ws.Range("K7").Value = 321.25
$ws.Range("L7").Value = 209.2
$ws.Range("M7").Value = -208.25
$ws.Range("N7").Value = -435.2
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
$ws.Range("H99").Value = 2038.2
$ws.Range("I99").Value = 2005.4546
$ws.Range("J99").Value = 2128.25
$ws.Range("K99").Value = 2005.4546
$ws.Range("L99").Value = 2128.25
$ws.Range("M99").Value = -507.4546
$ws.Range("N99").Value = -5124.25
$ws.Range("H126").Value = 2038.2
$ws.Range("I126").Value = 2005.4546
$ws.Range("J126").Value = 2128.25
$ws.Range("K126").Value = 6016.3638
$ws.Range("L126").Value = 6384.75
$ws.Range("M126").Value = -3546.3638
$ws.Range("N126").Value = -11324.75
$ws.Range("H133").Value = 46825.6
$ws.Range("J133").Value = 46825.6
$ws.Range("L133").Value = 46825.6
$ws.Range("N133").Value = -51885.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2911.25
$ws.Range("J75").Value = 2911.25
$ws.Range("L75").Value = 8733.75
$ws.Range("N75").Value = -10729.75
$ws.Range("H78").Value = 2911.25
$ws.Range("J78").Value = 2911.25
$ws.Range("L78").Value = 26201.25
$ws.Range("N78").Value = -36185.25
$ws.Range("H81").Value = 1866.2222
$ws.Range("J81").Value = 3010
$ws.Range("L81").Value = 9030
$ws.Range("N81").Value = -11276
$ws.Range("H84").Value = 1866.2222
$ws.Range("J84").Value = 3010
$ws.Range("L84").Value = 27090
$ws.Range("N84").Value = -38322
$ws.Range("H87").Value = 3010.6667
$ws.Range("I87").Value = 3000
$ws.Range("J87").Value = 3016
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 9048
$ws.Range("M87").Value = -7752
$ws.Range("N87").Value = -11544
$ws.Range("H90").Value = 3010.6667
$ws.Range("I90").Value = 3000
$ws.Range("J90").Value = 3016
$ws.Range("K90").Value = 27000
$ws.Range("L90").Value = 27144
$ws.Range("M90").Value = -20760
$ws.Range("N90").Value = -39624
$ws.Range("H98").Value = 725.7
$ws.Range("I98").Value = 101.6
$ws.Range("J98").Value = 1349.8
$ws.Range("K98").Value = 304.8
$ws.Range("L98").Value = 4049.4
$ws.Range("M98").Value = 1193.2
$ws.Range("N98").Value = -7045.4
$ws.Range("H131").Value = 47620996
$ws.Range("J131").Value = 2566.5
$ws.Range("L131").Value = 7699.5
$ws.Range("N131").Value = -17779.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17312098
$ws.Range("I70").Value = 22731204
$ws.Range("J70").Value = 13338087
$ws.Range("K70").Value = 22731204
$ws.Range("L70").Value = 13338087
$ws.Range("M70").Value = -22730934
$ws.Range("N70").Value = -13338627
$ws.Range("H73").Value = 17312098
$ws.Range("I73").Value = 22731204
$ws.Range("J73").Value = 13338087
$ws.Range("K73").Value = 22731204
$ws.Range("L73").Value = 13338087
$ws.Range("M73").Value = -22730268
$ws.Range("N73").Value = -13339959
$ws.Range("H104").Value = 48810.145
$ws.Range("J104").Value = 48810.145
$ws.Range("L104").Value = 48810.145
$ws.Range("N104").Value = -55798.145
$ws.Range("H122").Value = 2285.5881
$ws.Range("J122").Value = 3144.2222
$ws.Range("L122").Value = 9432.6666
$ws.Range("N122").Value = -14332.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1600.125
$ws.Range("J22").Value = 1685.8572
$ws.Range("L22").Value = 1685.8572
$ws.Range("N22").Value = -2275.8572
$ws.Range("H27").Value = 1600.125
$ws.Range("J27").Value = 1685.8572
$ws.Range("L27").Value = 1685.8572
$ws.Range("N27").Value = -1899.8572
$ws.Range("H40").Value = 2625.5715
$ws.Range("I40").Value = 2476
$ws.Range("K40").Value = 2476
$ws.Range("M40").Value = -2340
$ws.Range("H61").Value = 2528.375
$ws.Range("I61").Value = 2345.4
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 2345.4
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -2143.4
$ws.Range("N61").Value = -3237.3333
$ws.Range("H113").Value = 2528.375
$ws.Range("I113").Value = 2345.4
$ws.Range("J113").Value = 2833.3333
$ws.Range("K113").Value = 2345.4
$ws.Range("L113").Value = 2833.3333
$ws.Range("M113").Value = -175.4000000000001
$ws.Range("N113").Value = -7173.3333
$ws.Range("H132").Value = 50057.523
$ws.Range("I132").Value = 2127.2727
$ws.Range("K132").Value = 6381.8181
$ws.Range("M132").Value = -3851.8181
$ws.Range("H133").Value = 36616.168
$ws.Range("J133").Value = 36616.168
$ws.Range("L133").Value = 36616.168
$ws.Range("N133").Value = -41676.168
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10835110
$ws.Range("I122").Value = 13686002
$ws.Range("J122").Value = 1719.8
$ws.Range("K122").Value = 41058006
$ws.Range("L122").Value = 5159.4
$ws.Range("M122").Value = -41055556
$ws.Range("N122").Value = -10059.4
$ws.Range("H126").Value = 47620668
$ws.Range("J126").Value = 2023.4615
$ws.Range("L126").Value = 6070.3845
$ws.Range("N126").Value = -11010.3845
